# Update NATMI LR-pair sheet (Efna3-Epha3) with refreshed TPM-derived values.
# The sending/target cluster set changes from {ECs, FAPs, MuSCs} x {ECs, FAPs, MuSCs}
# (9 data rows) to {ECs, FAPs} x {ECs, FAPs, MuSCs, Resolving-Mac} (8 data rows),
# and every numeric column is refreshed with newly computed TPM-based statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old last data row (previously row 10, the MuSCs->MuSCs combination)
# so the sheet shrinks from 10 rows (incl. header) to 9.
$ws.Rows.Item(10).Delete()

# New data for rows 2-9 (20 columns, A:T) - cluster pair labels + refreshed stats.
$rowsData = @(
    @{ Row = 2; Values = @("ECs", "Efna3", "Epha3", "ECs", 3, 1, 0.3193606666666667, 0.958082, 0.7979421849584948, 0.7979421849584948, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.02420626109066666, 0.217856349816, 0.001420497439865801, 0.001420497439865801) },
    @{ Row = 3; Values = @("ECs", "Efna3", "Epha3", "FAPs", 3, 1, 0.3193606666666667, 0.958082, 0.7979421849584948, 0.7979421849584948, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 12.81509416041355, 115.335847443722, 0.7520289225305253, 0.7520289225305252) },
    @{ Row = 4; Values = @("ECs", "Efna3", "Epha3", "MuSCs", 3, 1, 0.3193606666666667, 0.958082, 0.7979421849584948, 0.7979421849584948, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 0.7486046095417778, 6.737441485875999, 0.04393040822549212, 0.04393040822549212) },
    @{ Row = 5; Values = @("ECs", "Efna3", "Epha3", "Resolving-Mac", 3, 1, 0.3193606666666667, 0.958082, 0.7979421849584948, 0.7979421849584948, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.009582949071111112, 0.08624654164000001, 0.0005623567626115689, 0.0005623567626115688) },
    @{ Row = 6; Values = @("FAPs", "Efna3", "Epha3", "ECs", 1, 0.3333333333333333, 0.08086966666666666, 0.242609, 0.2020578150415052, 0.2020578150415053, 2, 0.6666666666666666, 0.07579599999999999, 0.227388, 0.001780200955210419, 0.001780200955210419, 0.006129597254666665, 0.05516637529199999, 0.0003597035153446177, 0.0003597035153446178) },
    @{ Row = 7; Values = @("FAPs", "Efna3", "Epha3", "FAPs", 1, 0.3333333333333333, 0.08086966666666666, 0.242609, 0.2020578150415052, 0.2020578150415053, 3, 1, 40.12734033333333, 120.382021, 0.9424604146848589, 0.9424604146848587, 3.245084636976555, 29.205761732789, 0.1904314921543336, 0.1904314921543336) },
    @{ Row = 8; Values = @("FAPs", "Efna3", "Epha3", "MuSCs", 1, 0.3333333333333333, 0.08086966666666666, 0.242609, 0.2020578150415052, 0.2020578150415053, 3, 1, 2.344072666666666, 7.032217999999999, 0.05505462557763778, 0.05505462557763778, 0.1895643751957777, 1.706079376762, 0.01112421735214566, 0.01112421735214566) },
    @{ Row = 9; Values = @("FAPs", "Efna3", "Epha3", "Resolving-Mac", 1, 0.3333333333333333, 0.08086966666666666, 0.242609, 0.2020578150415052, 0.2020578150415053, 2, 0.6666666666666666, 0.03000666666666667, 0.09002, 0.0007047587822930054, 0.0007047587822930053, 0.002426629131111111, 0.02183966218, 0.0001424020196814365, 0.0001424020196814365) }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $vals = $rd.Values
    $arr = New-Object 'object[,]' 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    $ws.Range("A$r`:T$r").Value = $arr
}

Write-Host "Sheet updated: dimension now $($ws.UsedRange.Address())"
